# "Incremento clase Main en la planilla de métricas. Ejecución de Prueba."
#
# - Se agrega el incremento "Clase Main" (fila 21) a la tabla de
#   Desarrollo y correctivos, con sus tiempos y líneas de código.
# - Se carga el tiempo de corrección de errores lógicos del incremento
#   "Clase Carrera" (fila 20, columna L) que había quedado sin completar.
# - Se completa la sección "Ejecución de la Prueba" (fila 30) con sus
#   horarios de inicio y fin.
#
# El resto de la planilla (totales, resumen, gráfico) se recalcula solo
# a partir de estos valores porque son celdas con fórmulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Desarrollo y correctivos: incremento "Clase Carrera" (fila 20) ---
# Tiempo de Corrección E.L. que antes estaba en 0
$ws.Range("L20").Value = 0.0013888888888888889

# --- Desarrollo y correctivos: nuevo incremento "Clase Main" (fila 21) ---
$ws.Range("C21").Value = "Clase Main"
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 0.0013888888888888889
$ws.Range("H21").Value = 0.91388888888888886
$ws.Range("I21").Value = 0.91666666666666663
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 5

# --- Ejecución de la Prueba (fila 30) ---
$ws.Range("B30").Value = 0.0013888888888888889
$ws.Range("C30").Value = 0.91666666666666663
$ws.Range("D30").Value = 0.92083333333333339

# Selección que queda activa sobre la fila recién cargada
$null = $ws.Range("L21").Select()

$wb.Application.Calculate()
